$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (row 30) down to the new row 31
$ws.Range("A30:B30").Copy()
$ws.Range("A31:B31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new data values: date serial 46001 (12/10/2025) and error count 64
$ws.Range("A31").Value = 46001
$ws.Range("B31").Value = 64

# Match the selection state shown in the saved workbook
$ws.Range("A31:B31").Select()
